# Updates the cryptos list on Sheet1 (rows 2-51) to match the refreshed
# coinranking.com snapshot: Price (D) / Volume(1h) (E) updates throughout,
# plus a block of rows (27-31, 42-44, 46-51) whose Coin/Link/Price/Volume
# were reshuffled because the underlying ranking order changed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number -> hashtable of column letter -> new value.
$updates = [ordered]@{
    2 = @{ "D" = "41.722.09"; "E" = "  +5.53%  " }
    3 = @{ "D" = "2.236.26"; "E" = "  +3.27%  " }
    4 = @{ "D" = "1.01"; "E" = "  +0.63%  " }
    5 = @{ "D" = "226.73"; "E" = "  -0.84%  " }
    6 = @{ "D" = "0.617"; "E" = "  -2.64%  " }
    7 = @{ "D" = "61.39"; "E" = "  -3.53%  " }
    8 = @{ "E" = "  +0.22%  " }
    9 = @{ "D" = "0.400"; "E" = "  +1.18%  " }
    10 = @{ "D" = "62.41"; "E" = "  +7.47%  " }
    11 = @{ "D" = "0.0877"; "E" = "  +2.91%  " }
    12 = @{ "D" = "0.104"; "E" = "  +0.07%  " }
    13 = @{ "D" = "2.565.25"; "E" = "  +3.07%  " }
    14 = @{ "D" = "15.55"; "E" = "  -3.27%  " }
    15 = @{ "D" = "21.78"; "E" = "  -1.19%  " }
    16 = @{ "D" = "0.796"; "E" = "  -1.81%  " }
    17 = @{ "D" = "5.53"; "E" = "  +0.34%  " }
    18 = @{ "D" = "2.216.68"; "E" = "  +2.14%  " }
    19 = @{ "D" = "41.477.85"; "E" = "  +5.03%  " }
    20 = @{ "D" = "73.07"; "E" = "  +1.33%  " }
    21 = @{ "D" = "0.0₃0893"; "E" = "  +5.40%  " }
    22 = @{ "D" = "5.98"; "E" = "  -3.89%  " }
    23 = @{ "D" = "247.08"; "E" = "  +7.40%  " }
    24 = @{ "D" = "0.999"; "E" = "  -0.14%  " }
    25 = @{ "D" = "2.40"; "E" = "  +2.67%  " }
    26 = @{ "D" = "2.33"; "E" = "  -0.85%  " }
    27 = @{ "B" = "Cosmos"; "C" = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"; "D" = "9.50"; "E" = "  -0.38%  " }
    28 = @{ "B" = "Monero"; "C" = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"; "D" = "168.47"; "E" = "  -2.10%  " }
    29 = @{ "B" = "Kaspa"; "C" = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"; "D" = "0.141"; "E" = "  +1.74%  " }
    30 = @{ "B" = "EthereumClassic"; "C" = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"; "D" = "20.02"; "E" = "  +0.67%  " }
    31 = @{ "B" = "ImmutableX"; "C" = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"; "D" = "1.41"; "E" = "  -0.88%  " }
    32 = @{ "D" = "2.83"; "E" = "  +6.36%  " }
    33 = @{ "D" = "0.121"; "E" = "  -1.07%  " }
    34 = @{ "D" = "4.91"; "E" = "  +4.24%  " }
    35 = @{ "D" = "4.59"; "E" = "  -0.20%  " }
    36 = @{ "D" = "0.0621"; "E" = "  +0.13%  " }
    37 = @{ "D" = "6.62"; "E" = "  -5.97%  " }
    38 = @{ "D" = "3.68"; "E" = "  -0.44%  " }
    39 = @{ "D" = "2.36"; "E" = "  -3.12%  " }
    40 = @{ "D" = "1.02"; "E" = "  +1.37%  " }
    41 = @{ "D" = "4.88"; "E" = "  +8.70%  " }
    42 = @{ "B" = "VeChain"; "C" = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"; "D" = "0.0234"; "E" = "  +2.67%  " }
    43 = @{ "B" = "FraxShare"; "C" = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"; "D" = "8.47"; "E" = "  +9.14%  " }
    44 = @{ "B" = "TerraClassic"; "C" = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"; "D" = "0.000221"; "E" = "  +22.56%  " }
    45 = @{ "D" = "99.28"; "E" = "  -3.27%  " }
    46 = @{ "B" = "Cronos"; "C" = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"; "D" = "0.0964"; "E" = "  +4.29%  " }
    47 = @{ "B" = "Maker"; "C" = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"; "D" = "1.479.43"; "E" = "  -2.91%  " }
    48 = @{ "B" = "HuobiToken"; "C" = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"; "D" = "2.82"; "E" = "  +0.40%  " }
    49 = @{ "B" = "TrustWalletToken"; "C" = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"; "D" = "1.16"; "E" = "  -3.92%  " }
    50 = @{ "B" = "InjectiveProtocol"; "C" = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"; "D" = "16.35"; "E" = "  -8.93%  " }
    51 = @{ "B" = "ARBITRUM"; "C" = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"; "D" = "1.06"; "E" = "  -4.08%  " }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $cellRef = "$col$row"
        $newValue = $cols[$col]
        $cell = $ws.Range($cellRef)
        if ($col -eq "D" -and $newValue -match "^[+-]?\d+(\.\d+)?$") {
            # Plain decimal text (e.g. "1.01") would otherwise be auto-coerced
            # to a number by the Value setter. Force text storage, matching the
            # source inlineStr cell, then drop the temporary "@" format so the
            # cell keeps its original (unstyled) appearance.
            $cell.NumberFormat = "@"
            $cell.Value = $newValue
            $cell.ClearFormats()
        } else {
            $cell.Value = $newValue
        }
    }
}
